$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colValues = @(0.8149060806616717,0.7750184406520191,0.7508280886198406,0.7410466867682999,0.7394271257576008,0.7506958630696374,0.8010908048804595,0.9022759107380409,0.9780288482903075,1.01279224614251,1.025999189961169,1.023152947089272,1.013877936018474,1.008202274758389,0.9757630097248295,0.9559396290366067,0.9445663113765477,0.9407204313138209,0.9580469108108787,1.016601076190142,1.055118717866435,1.034538589335796,0.9570941357488891,0.8746527075001893)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colValues[$i]
}

$colValues = @(0.1933857673950889,0.1939736074560976,0.1943577008368997,0.1945200611728879,0.194547374130277,0.1943598668212552,0.193583658043055,0.1922445352046722,0.1913712635613543,0.1909977923007524,0.190859772546478,0.1908893463300174,0.1909863691602354,0.1910462415424661,0.1913961481501438,0.1916168859373961,0.191746088029177,0.1917902187258811,0.1915931563427336,0.1909577788882437,0.1905623674476402,0.190771594945236,0.1916038773366964,0.1925873129601907)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colValues[$i]
}

$colValues = @(0.2281357363757017,0.2268825072741194,0.2261996318757404,0.2259431998857586,0.2259019408881606,0.2261960850057605,0.227685685295544,0.2312914691581369,0.2343551271377819,0.2358382808251918,0.2364127201518684,0.2362884360781692,0.235885284116911,0.2356400070375315,0.2342599944928949,0.2334362726444503,0.2329709158804398,0.232814803302098,0.2335230876332162,0.236003352725092,0.2376989339698596,0.2367871676533397,0.2334838129930219,0.230243010044866)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colValues[$i]
}

$colValues = @(2.202178076897397,2.200917907954903,2.201169440998598,2.201530074588845,2.201605557364871,2.201173259044054,2.201530957491755,2.210356230186804,2.221781246709696,2.228049526615308,2.230576913547253,2.230025763050222,2.228254376983926,2.227189363110242,2.221393126812771,2.218111428487461,2.216324705178096,2.215737076314255,2.218450338128022,2.228770506188454,2.23641123095058,2.23225135050393,2.218296805663343,2.207100226966006)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $colValues[$i]
}

$colValues = @(0.002484787234215415,0.002487143628304584,0.002488668306307481,0.002489309256433648,0.002489416873071518,0.002488676870789281,0.002485583597191592,0.002480132723350525,0.002476499230776338,0.002474926096785909,0.002474341802290825,0.00247446713361554,0.002474877797861478,0.00247513082721223,0.002476603639275858,0.002477527553218799,0.002478066474433481,0.002478250235469917,0.00247742842397643,0.002474756866363172,0.002473077373843946,0.002473967680848168,0.002477473216240808,0.002481541867647916)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colValues[$i]
}

$colValues = @(0.06720758584429021,0.06638006353365711,0.06586880983464738,0.06565968885430706,0.06562491764685419,0.06586599270258375,0.06692291700088049,0.0689701154079394,0.07045825642043368,0.0711317040008872,0.07138620505511284,0.07133141699402756,0.07115265243341184,0.07104308599794606,0.07041417315435794,0.07002744562093355,0.06980468011067487,0.06972919925626542,0.07006864768403531,0.07120517410949034,0.07194492881797032,0.07155038960786086,0.07005002157199058,0.06841906076351378)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $colValues[$i]
}

$colValues = @(0.3578018687996689,0.3253221577429599,0.3054464656298137,0.2973642784447748,0.2960232974374151,0.3053373956734333,0.3465892373754969,0.4279986314517146,0.4881072860000302,0.5155136403542429,0.5259003151301727,0.5236629899021636,0.5163679913159456,0.5119006832987907,0.4863174318327879,0.4706386040666928,0.4616264829022327,0.4585761672156821,0.4723070325257481,0.5185104831669207,0.5487563130542412,0.532609214870547,0.4715527302305986,0.4059218494859351)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $colValues[$i]
}

$colValues = @(0.4304720508122273,0.4213189132172204,0.4158751048123719,0.4137012343948498,0.4133429601391612,0.4158456066805485,0.4272795865308865,0.4510925871223037,0.4694281156791078,0.4779503060856314,0.4812033252748051,0.480501583106232,0.4782174170806144,0.4768216595573307,0.4688747996467697,0.4640459321477266,0.4612855698375995,0.4603538977128636,0.4645582077430817,0.4788876321502897,0.488403316487009,0.4833109100789699,0.4643265586607868,0.4445024361192225)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $colValues[$i]
}

$colValues = @(2.001538814254022,2.020873338268105,2.033348906436248,2.038584740028922,2.039463323959224,2.033418903481598,2.008080001178174,1.963181403777261,1.933110848054159,1.920063355910628,1.915213428364861,1.916253904873801,1.919662528014355,1.921762243572321,1.93397624727986,1.941630970138325,1.946093226055433,1.947614281023423,1.940809957547208,1.918658866215897,1.904711549579827,1.912107016824955,1.941180945913957,1.974815212490963)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $colValues[$i]
}

$colValues = @(4.088716896165778,4.111959272972399,4.128099907232468,4.135147716158087,4.136346414924304,4.128193051660716,4.09634293152935,4.048712725401003,4.02275057394499,4.012899504219945,4.009450759896964,4.010180984087867,4.012610129963122,4.01413472725676,4.023433781372347,4.029640191268669,4.033394372910436,4.03469715133383,4.028960422369238,4.011888988003847,4.00237340057086,4.007301874876532,4.029267166392884,4.060011339560759)
for ($i = 0; $i -lt 24; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $colValues[$i]
}
